$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 103.4275383333333
$ws.Cells.Item(2, 8).Value = 310.282615
$ws.Cells.Item(2, 9).Value = 0.2485530285127421
$ws.Cells.Item(2, 10).Value = 0.2485530285127421
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 20.02757633333333
$ws.Cells.Item(2, 14).Value = 60.082729
$ws.Cells.Item(2, 15).Value = 0.200201311135073
$ws.Cells.Item(2, 16).Value = 0.200201311135073
$ws.Cells.Item(2, 17).Value = 2071.402918939593
$ws.Cells.Item(2, 18).Value = 18642.62627045634
$ws.Cells.Item(2, 19).Value = 0.04976064219484416
$ws.Cells.Item(2, 20).Value = 0.04976064219484415
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 103.4275383333333
$ws.Cells.Item(3, 8).Value = 310.282615
$ws.Cells.Item(3, 9).Value = 0.2485530285127421
$ws.Cells.Item(3, 10).Value = 0.2485530285127421
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 33.29907966666667
$ws.Cells.Item(3, 14).Value = 99.89723900000001
$ws.Cells.Item(3, 15).Value = 0.3328670078646686
$ws.Cells.Item(3, 16).Value = 0.3328670078646686
$ws.Cells.Item(3, 17).Value = 3444.041838688888
$ws.Cells.Item(3, 18).Value = 30996.37654819999
$ws.Cells.Item(3, 19).Value = 0.08273510289673813
$ws.Cells.Item(3, 20).Value = 0.08273510289673811
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 103.4275383333333
$ws.Cells.Item(4, 8).Value = 310.282615
$ws.Cells.Item(4, 9).Value = 0.2485530285127421
$ws.Cells.Item(4, 10).Value = 0.2485530285127421
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 33.93321599999999
$ws.Cells.Item(4, 14).Value = 101.799648
$ws.Cells.Item(4, 15).Value = 0.3392060138062123
$ws.Cells.Item(4, 16).Value = 0.3392060138062122
$ws.Cells.Item(4, 17).Value = 3509.62899861328
$ws.Cells.Item(4, 18).Value = 31586.66098751952
$ws.Cells.Item(4, 19).Value = 0.08431068202126908
$ws.Cells.Item(4, 20).Value = 0.08431068202126905
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 103.4275383333333
$ws.Cells.Item(5, 8).Value = 310.282615
$ws.Cells.Item(5, 9).Value = 0.2485530285127421
$ws.Cells.Item(5, 10).Value = 0.2485530285127421
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 12.77731666666667
$ws.Cells.Item(5, 14).Value = 38.33195
$ws.Cells.Item(5, 15).Value = 0.1277256671940461
$ws.Cells.Item(5, 16).Value = 0.1277256671940461
$ws.Cells.Item(5, 17).Value = 1321.526409338806
$ws.Cells.Item(5, 18).Value = 11893.73768404925
$ws.Cells.Item(5, 19).Value = 0.03174660139989075
$ws.Cells.Item(5, 20).Value = 0.03174660139989074
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 216.130539
$ws.Cells.Item(6, 8).Value = 648.391617
$ws.Cells.Item(6, 9).Value = 0.5193964865470273
$ws.Cells.Item(6, 10).Value = 0.5193964865470272
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 20.02757633333333
$ws.Cells.Item(6, 14).Value = 60.082729
$ws.Cells.Item(6, 15).Value = 0.200201311135073
$ws.Cells.Item(6, 16).Value = 0.200201311135073
$ws.Cells.Item(6, 17).Value = 4328.570867786977
$ws.Cells.Item(6, 18).Value = 38957.13781008279
$ws.Cells.Item(6, 19).Value = 0.1039838576056652
$ws.Cells.Item(6, 20).Value = 0.1039838576056652
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 216.130539
$ws.Cells.Item(7, 8).Value = 648.391617
$ws.Cells.Item(7, 9).Value = 0.5193964865470273
$ws.Cells.Item(7, 10).Value = 0.5193964865470272
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 33.29907966666667
$ws.Cells.Item(7, 14).Value = 99.89723900000001
$ws.Cells.Item(7, 15).Value = 0.3328670078646686
$ws.Cells.Item(7, 16).Value = 0.3328670078646686
$ws.Cells.Item(7, 17).Value = 7196.948036560608
$ws.Cells.Item(7, 18).Value = 64772.53232904547
$ws.Cells.Item(7, 19).Value = 0.1728899543723306
$ws.Cells.Item(7, 20).Value = 0.1728899543723305
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 216.130539
$ws.Cells.Item(8, 8).Value = 648.391617
$ws.Cells.Item(8, 9).Value = 0.5193964865470273
$ws.Cells.Item(8, 10).Value = 0.5193964865470272
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 33.93321599999999
$ws.Cells.Item(8, 14).Value = 101.799648
$ws.Cells.Item(8, 15).Value = 0.3392060138062123
$ws.Cells.Item(8, 16).Value = 0.3392060138062122
$ws.Cells.Item(8, 17).Value = 7334.004264083423
$ws.Cells.Item(8, 18).Value = 66006.0383767508
$ws.Cells.Item(8, 19).Value = 0.1761824117865691
$ws.Cells.Item(8, 20).Value = 0.176182411786569
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 216.130539
$ws.Cells.Item(9, 8).Value = 648.391617
$ws.Cells.Item(9, 9).Value = 0.5193964865470273
$ws.Cells.Item(9, 10).Value = 0.5193964865470272
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 12.77731666666667
$ws.Cells.Item(9, 14).Value = 38.33195
$ws.Cells.Item(9, 15).Value = 0.1277256671940461
$ws.Cells.Item(9, 16).Value = 0.1277256671940461
$ws.Cells.Item(9, 17).Value = 2761.56833814035
$ws.Cells.Item(9, 18).Value = 24854.11504326315
$ws.Cells.Item(9, 19).Value = 0.06634026278246245
$ws.Cells.Item(9, 20).Value = 0.06634026278246243
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 71.607325
$ws.Cells.Item(10, 8).Value = 214.821975
$ws.Cells.Item(10, 9).Value = 0.1720839321833696
$ws.Cells.Item(10, 10).Value = 0.1720839321833696
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 20.02757633333333
$ws.Cells.Item(10, 14).Value = 60.082729
$ws.Cells.Item(10, 15).Value = 0.200201311135073
$ws.Cells.Item(10, 16).Value = 0.200201311135073
$ws.Cells.Item(10, 17).Value = 1434.121167463308
$ws.Cells.Item(10, 18).Value = 12907.09050716978
$ws.Cells.Item(10, 19).Value = 0.03445142884838957
$ws.Cells.Item(10, 20).Value = 0.03445142884838957
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 71.607325
$ws.Cells.Item(11, 8).Value = 214.821975
$ws.Cells.Item(11, 9).Value = 0.1720839321833696
$ws.Cells.Item(11, 10).Value = 0.1720839321833696
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 33.29907966666667
$ws.Cells.Item(11, 14).Value = 99.89723900000001
$ws.Cells.Item(11, 15).Value = 0.3328670078646686
$ws.Cells.Item(11, 16).Value = 0.3328670078646686
$ws.Cells.Item(11, 17).Value = 2384.458019891892
$ws.Cells.Item(11, 18).Value = 21460.12217902703
$ws.Cells.Item(11, 19).Value = 0.05728106360746477
$ws.Cells.Item(11, 20).Value = 0.05728106360746477
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 71.607325
$ws.Cells.Item(12, 8).Value = 214.821975
$ws.Cells.Item(12, 9).Value = 0.1720839321833696
$ws.Cells.Item(12, 10).Value = 0.1720839321833696
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 33.93321599999999
$ws.Cells.Item(12, 14).Value = 101.799648
$ws.Cells.Item(12, 15).Value = 0.3392060138062123
$ws.Cells.Item(12, 16).Value = 0.3392060138062122
$ws.Cells.Item(12, 17).Value = 2429.8668264072
$ws.Cells.Item(12, 18).Value = 21868.8014376648
$ws.Cells.Item(12, 19).Value = 0.05837190467601935
$ws.Cells.Item(12, 20).Value = 0.05837190467601935
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 71.607325
$ws.Cells.Item(13, 8).Value = 214.821975
$ws.Cells.Item(13, 9).Value = 0.1720839321833696
$ws.Cells.Item(13, 10).Value = 0.1720839321833696
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 12.77731666666667
$ws.Cells.Item(13, 14).Value = 38.33195
$ws.Cells.Item(13, 15).Value = 0.1277256671940461
$ws.Cells.Item(13, 16).Value = 0.1277256671940461
$ws.Cells.Item(13, 17).Value = 914.9494671779166
$ws.Cells.Item(13, 18).Value = 8234.545204601251
$ws.Cells.Item(13, 19).Value = 0.02197953505149586
$ws.Cells.Item(13, 20).Value = 0.02197953505149585
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 24.953198
$ws.Cells.Item(14, 8).Value = 74.859594
$ws.Cells.Item(14, 9).Value = 0.05996655275686102
$ws.Cells.Item(14, 10).Value = 0.05996655275686102
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 20.02757633333333
$ws.Cells.Item(14, 14).Value = 60.082729
$ws.Cells.Item(14, 15).Value = 0.200201311135073
$ws.Cells.Item(14, 16).Value = 0.200201311135073
$ws.Cells.Item(14, 17).Value = 499.7520777057807
$ws.Cells.Item(14, 18).Value = 4497.768699352026
$ws.Cells.Item(14, 19).Value = 0.0120053824861741
$ws.Cells.Item(14, 20).Value = 0.0120053824861741
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 24.953198
$ws.Cells.Item(15, 8).Value = 74.859594
$ws.Cells.Item(15, 9).Value = 0.05996655275686102
$ws.Cells.Item(15, 10).Value = 0.05996655275686102
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 33.29907966666667
$ws.Cells.Item(15, 14).Value = 99.89723900000001
$ws.Cells.Item(15, 15).Value = 0.3328670078646686
$ws.Cells.Item(15, 16).Value = 0.3328670078646686
$ws.Cells.Item(15, 17).Value = 830.9185281401075
$ws.Cells.Item(15, 18).Value = 7478.266753260968
$ws.Cells.Item(15, 19).Value = 0.01996088698813512
$ws.Cells.Item(15, 20).Value = 0.01996088698813512
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 24.953198
$ws.Cells.Item(16, 8).Value = 74.859594
$ws.Cells.Item(16, 9).Value = 0.05996655275686102
$ws.Cells.Item(16, 10).Value = 0.05996655275686102
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 33.93321599999999
$ws.Cells.Item(16, 14).Value = 101.799648
$ws.Cells.Item(16, 15).Value = 0.3392060138062123
$ws.Cells.Item(16, 16).Value = 0.3392060138062122
$ws.Cells.Item(16, 17).Value = 846.7422576247678
$ws.Cells.Item(16, 18).Value = 7620.680318622912
$ws.Cells.Item(16, 19).Value = 0.02034101532235476
$ws.Cells.Item(16, 20).Value = 0.02034101532235475
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 24.953198
$ws.Cells.Item(17, 8).Value = 74.859594
$ws.Cells.Item(17, 9).Value = 0.05996655275686102
$ws.Cells.Item(17, 10).Value = 0.05996655275686102
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 12.77731666666667
$ws.Cells.Item(17, 14).Value = 38.33195
$ws.Cells.Item(17, 15).Value = 0.1277256671940461
$ws.Cells.Item(17, 16).Value = 0.1277256671940461
$ws.Cells.Item(17, 17).Value = 318.8349126920333
$ws.Cells.Item(17, 18).Value = 2869.5142142283
$ws.Cells.Item(17, 19).Value = 0.007659267960197039
$ws.Cells.Item(17, 20).Value = 0.007659267960197036
